# Apply Natmi Efna1-Epha2 LR-pair updates per Dr Hou advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.905615
$ws.Range("H2").Value = 35.716845
$ws.Range("I2").Value = 0.8197078149061106
$ws.Range("J2").Value = 0.8197078149061106
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 11.55727433333333
$ws.Range("N2").Value = 34.671823
$ws.Range("O2").Value = 0.5239815261112396
$ws.Range("P2").Value = 0.5239815261112395
$ws.Range("Q2").Value = 137.5964586620483
$ws.Range("R2").Value = 1238.368127958435
$ws.Range("S2").Value = 0.4295117518198133
$ws.Range("T2").Value = 0.4295117518198132
# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.905615
$ws.Range("H3").Value = 35.716845
$ws.Range("I3").Value = 0.8197078149061106
$ws.Range("J3").Value = 0.8197078149061106
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.09477133333333332
$ws.Range("N3").Value = 0.284314
$ws.Range("O3").Value = 0.004296724853919303
$ws.Range("P3").Value = 0.004296724853919302
$ws.Range("Q3").Value = 1.128311007703333
$ws.Range("R3").Value = 10.15479906933
$ws.Range("S3").Value = 0.003522058941258969
$ws.Range("T3").Value = 0.003522058941258969
# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.905615
$ws.Range("H4").Value = 35.716845
$ws.Range("I4").Value = 0.8197078149061106
$ws.Range("J4").Value = 0.8197078149061106
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.404599
$ws.Range("N4").Value = 31.213797
$ws.Range("O4").Value = 0.471721749034841
$ws.Range("P4").Value = 0.471721749034841
$ws.Range("Q4").Value = 123.873149923385
$ws.Range("R4").Value = 1114.858349310465
$ws.Range("S4").Value = 0.3866740041450382
$ws.Range("T4").Value = 0.3866740041450382
# Row 5
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.63134
$ws.Range("H5").Value = 4.894019999999999
$ws.Range("I5").Value = 0.1123186115768849
$ws.Range("J5").Value = 0.1123186115768849
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 11.55727433333333
$ws.Range("N5").Value = 34.671823
$ws.Range("O5").Value = 0.5239815261112396
$ws.Range("P5").Value = 0.5239815261112395
$ws.Range("Q5").Value = 18.85384391094
$ws.Range("R5").Value = 169.68459519846
$ws.Range("S5").Value = 0.05885287750475169
$ws.Range("T5").Value = 0.05885287750475168
# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.63134
$ws.Range("H6").Value = 4.894019999999999
$ws.Range("I6").Value = 0.1123186115768849
$ws.Range("J6").Value = 0.1123186115768849
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.09477133333333332
$ws.Range("N6").Value = 0.284314
$ws.Range("O6").Value = 0.004296724853919303
$ws.Range("P6").Value = 0.004296724853919302
$ws.Range("Q6").Value = 0.15460426692
$ws.Range("R6").Value = 1.39143840228
$ws.Range("S6").Value = 0.0004826021699201097
$ws.Range("T6").Value = 0.0004826021699201096
# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.63134
$ws.Range("H7").Value = 4.894019999999999
$ws.Range("I7").Value = 0.1123186115768849
$ws.Range("J7").Value = 0.1123186115768849
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.404599
$ws.Range("N7").Value = 31.213797
$ws.Range("O7").Value = 0.471721749034841
$ws.Range("P7").Value = 0.471721749034841
$ws.Range("Q7").Value = 16.97343853266
$ws.Range("R7").Value = 152.76094679394
$ws.Range("S7").Value = 0.05298313190221308
$ws.Range("T7").Value = 0.05298313190221308
# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efna1"
$ws.Range("C8").Value = "Epha2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.9872630000000001
$ws.Range("H8").Value = 2.961789
$ws.Range("I8").Value = 0.0679735735170045
$ws.Range("J8").Value = 0.0679735735170045
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 11.55727433333333
$ws.Range("N8").Value = 34.671823
$ws.Range("O8").Value = 0.5239815261112396
$ws.Range("P8").Value = 0.5239815261112395
$ws.Range("Q8").Value = 11.41006933014967
$ws.Range("R8").Value = 102.690623971347
$ws.Range("S8").Value = 0.03561689678667456
$ws.Range("T8").Value = 0.03561689678667455
# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efna1"
$ws.Range("C9").Value = "Epha2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.9872630000000001
$ws.Range("H9").Value = 2.961789
$ws.Range("I9").Value = 0.0679735735170045
$ws.Range("J9").Value = 0.0679735735170045
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09477133333333332
$ws.Range("N9").Value = 0.284314
$ws.Range("O9").Value = 0.004296724853919303
$ws.Range("P9").Value = 0.004296724853919302
$ws.Range("Q9").Value = 0.09356423086066666
$ws.Range("R9").Value = 0.842078077746
$ws.Range("S9").Value = 0.0002920637427402242
$ws.Range("T9").Value = 0.0002920637427402242
# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efna1"
$ws.Range("C10").Value = "Epha2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.9872630000000001
$ws.Range("H10").Value = 2.961789
$ws.Range("I10").Value = 0.0679735735170045
$ws.Range("J10").Value = 0.0679735735170045
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 10.404599
$ws.Range("N10").Value = 31.213797
$ws.Range("O10").Value = 0.471721749034841
$ws.Range("P10").Value = 0.471721749034841
$ws.Range("Q10").Value = 10.272075622537
$ws.Range("R10").Value = 92.44868060283301
$ws.Range("S10").Value = 0.03206461298758972
$ws.Range("T10").Value = 0.03206461298758972
